$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to insert before the current row 2 (9 rows)
$newTopRows = @(
    @(0.0189368221908807, -0.030695978552103, -0.0142026171088218),
    @(-0.0242818929255008, 0.0024434609804302, 0.0209221355617046),
    @(-0.0287106670439243, -0.0282525178045034, 0.00534507073462),
    @(0.0103847095742821, -0.0001527163112768, 0.0277943685650825),
    @(-0.0300851128995418, -0.0401643887162208, 0.0236710291355848),
    @(-0.0271835029125213, -0.0343611687421798, 0.0047342055477201),
    @(-0.0245873257517814, -0.0226020142436027, -0.0163406450301408),
    @(-0.0479529201984405, 0.07635815441608421, -0.1252273768186569),
    @(0.1218676194548606, 0.3381139039993286, -0.0650571510195732)
)

# Insert 9 new rows at row 2, shifting existing data down
$insertRange = $ws.Range("2:10")
$insertRange.Insert(-4121)  # xlShiftDown = -4121
$ws.Range("A2:C10").ClearFormats()

for ($i = 0; $i -lt $newTopRows.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newTopRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newTopRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newTopRows[$i][2]
}

# Append a new row at the end (row 31)
$ws.Cells.Item(31, 1).Value = -0.0186313893646001
$ws.Cells.Item(31, 2).Value = -0.1244637966156005
$ws.Cells.Item(31, 3).Value = 0.0003054326225537
